$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row before sheet row 64 (a Table1 data row). This shifts
#     the existing rows 64-135 down to 65-136 (values/formulas travel with
#     them), matching the "Leave Card" gaining a 2024 section.
$ws.Rows("64:64").Insert()

# The freshly-inserted row 64 comes back with generic placeholder styles;
# restore the normal data-row look by copying formats from row 65 (which
# now holds what used to be row 64's formatting/content).
$ws.Range("A65:K65").Copy()
$ws.Range("A64:K64").PasteSpecial(-4122)

# Row 64 becomes a new year-header row, just like row 49 ("2023"). Copy
# that row's PERIOD-column styling (centered, quote-prefixed text format)
# then write the "2024" label and restore the EARNED helper formula.
$ws.Range("A49").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A64").Value = "'2024"
$ws.Range("G64").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Row 65 (previously row 64, dated 45292) gains a Sick Leave entry that
# spans the 2023/2024 year boundary.
$ws.Range("B65").Value = "SL(3-0-0)"
$ws.Range("H65").Value = 3
$ws.Range("K65").Value = "12/30/2023 - 1/1/2024"

# EARNED amounts posted for three prior periods.
$ws.Range("C61").Value = 1.25
$ws.Range("C62").Value = 1.25
$ws.Range("C63").Value = 1.25

# Grow Table1 so it covers the newly-inserted row.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K136"))

# The resize leaves the last row's calculated-column formula in the
# abbreviated [@EARNED] form; re-apply it explicitly so it serializes the
# same way as the rest of the column.
$ws.Range("G136").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Match the author's final on-screen selection.
[void]$ws.Range("K65").Select()
